$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F11").Value = 2086
$wsExhibit.Range("F14").Value = 1340
$wsExhibit.Range("F22").Value = 55
$wsExhibit.Range("F25").Value = 1126

# Sheet "全部类型" updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F12").Value = 2086
$wsAll.Range("F15").Value = 1340
$wsAll.Range("F23").Value = 55
$wsAll.Range("F26").Value = 1126
